# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets.
# All values are written as literal text (matching the existing inline-string
# data already in the sheets), so date-, time- and percentage-looking values
# must be prevented from being auto-converted into Excel dates/numbers.

$wb = $excel.ActiveWorkbook

# ---- PIR sheet: append rows 78-85 ----
$ws = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @("2026-02-06","09:46:24","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:46:29","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:46:34","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:46:35","09:00","Bathroom","Motion Detected","Active"),
    @("2026-02-06","09:46:43","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:46:48","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:46:53","09:00","Bathroom","No Motion","Inactive"),
    @("2026-02-06","09:46:57","09:00","Bathroom","Motion Detected","Active")
)
$startRow = 78
for ($i = 0; $i -lt $pirRows.Length; $i++) {
    $rowvals = $pirRows[$i]
    $r = $startRow + $i
    for ($col = 1; $col -le $rowvals.Length; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $rowvals[$col - 1]
        $cell.ClearFormats()
    }
}

# ---- Humidity sheet: append rows 20-25 ----
$ws = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @("2026-02-06","09:46:19","09:00","Bathroom","73.2%","Active"),
    @("2026-02-06","09:46:24","09:00","Bathroom","72.8%","Active"),
    @("2026-02-06","09:46:39","09:00","Bathroom","93.6%","Active"),
    @("2026-02-06","09:46:44","09:00","Bathroom","89.3%","Active"),
    @("2026-02-06","09:46:49","09:00","Bathroom","82.0%","Active"),
    @("2026-02-06","09:46:54","09:00","Bathroom","84.5%","Active")
)
$startRow = 20
for ($i = 0; $i -lt $humidityRows.Length; $i++) {
    $rowvals = $humidityRows[$i]
    $r = $startRow + $i
    for ($col = 1; $col -le $rowvals.Length; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $rowvals[$col - 1]
        $cell.ClearFormats()
    }
}

# ---- Temperature sheet: append rows 20-25 ----
$ws = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @("2026-02-06","09:46:19","09:00","Bathroom","27.3C","Active"),
    @("2026-02-06","09:46:24","09:00","Bathroom","27.3C","Active"),
    @("2026-02-06","09:46:39","09:00","Bathroom","27.4C","Active"),
    @("2026-02-06","09:46:44","09:00","Bathroom","27.4C","Active"),
    @("2026-02-06","09:46:49","09:00","Bathroom","27.5C","Active"),
    @("2026-02-06","09:46:54","09:00","Bathroom","27.4C","Active")
)
$startRow = 20
for ($i = 0; $i -lt $temperatureRows.Length; $i++) {
    $rowvals = $temperatureRows[$i]
    $r = $startRow + $i
    for ($col = 1; $col -le $rowvals.Length; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $rowvals[$col - 1]
        $cell.ClearFormats()
    }
}
